# "update to new format inv"
# The laptop-inventory numbers in column C (inventory_number) are being
# migrated from the old "PPAHONB0xx" scheme to the new "HO-NB-COE-0xx"
# scheme, and the "nrp pengguna" (user NRP) column U is updated to the
# new NRP 22003193 for all three sample rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column C - inventory_number: old PPAHONB0xx -> new HO-NB-COE-0xx format
$ws.Range("C25").Value = "HO-NB-COE-001"
$ws.Range("C26").Value = "HO-NB-COE-002"
$ws.Range("C27").Value = "HO-NB-COE-003"

# Column U - nrp pengguna: updated to the new NRP for all three rows
$ws.Range("U25").Value = 22003193
$ws.Range("U26").Value = 22003193
$ws.Range("U27").Value = 22003193

# Reflect the author's last on-screen selection/scroll position before save
$ws.Range("U27").Select()
